# roberta's robustness test results
#
# The sheet holds token-level NER predictions grouped by `message_id`
# (column A), where each message occupies a contiguous block of rows.
# This edit reorders the four message blocks and renumbers message_id
# sequentially (0..3) to match the new block order, while every other
# column (message text, token, token_index, labels, booleans, entity
# types, error_type) simply travels along with its original row.
#
# New block order (by original message_id): Compass(2), GPS(3),
# High wind(0), Motor speed(1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 12
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# --- Read all data rows (everything below the header row) ---
$data = @()
for ($r = 2; $r -le $lastRow; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowVals += ,($ws.Cells.Item($r, $c).Value2)
    }
    $data += ,$rowVals
}

# --- Group consecutive rows into per-message blocks using column A ---
$blocks = @()
$curBlock = @()
$curId = $null
foreach ($row in $data) {
    $mid = $row[0]
    if ($curId -ne $null -and $mid -ne $curId) {
        $blocks += ,$curBlock
        $curBlock = @()
    }
    $curBlock += ,$row
    $curId = $mid
}
if ($curBlock.Count -gt 0) {
    $blocks += ,$curBlock
}

# --- Reorder the blocks and renumber message_id sequentially ---
$newOrder = @(2, 3, 0, 1)

$newData = @()
for ($i = 0; $i -lt $newOrder.Count; $i++) {
    $blk = $blocks[$newOrder[$i]]
    foreach ($row in $blk) {
        $newRow = @($i, $row[1], $row[2], $row[3], $row[4], $row[5], $row[6], $row[7], $row[8], $row[9], $row[10], $row[11])
        $newData += ,$newRow
    }
}

# --- Write the reordered data back over the original range ---
for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $i + 2
    $row = $newData[$i]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
